# Arreglos en los scripts de anulacion y rehabilitacion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: keep environment (i-preproducciongestion / su / silverarrow) but
# update the policy number. Use a leading apostrophe so it is stored as
# text (preserving the leading zero) instead of being coerced to a number.
$ws.Range("E2").Value = "'04104013566"

# Row 3 previously held the "ssurgwsoadev4" test environment; it is being
# retired. Clear all of its cell contents. B3/E3 keep their existing
# styles (hyperlink-like font / grey font respectively) even once empty.
$ws.Range("A3").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()

# Row 5 previously held the "i-gestion-ssur-oci" environment (with a
# hyperlink on B5); it is being retired too. Clear its contents, keeping
# B5's existing style even once empty.
$ws.Range("A5").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()

# Remove the now-orphaned hyperlink that was attached to B5.
$ws.Hyperlinks.Delete()

# Move the active selection from the old E3 to the now-current E2.
[void]$ws.Range("E2").Select()
